$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the wording of the "Определение" (Definition) cell B12: the old text
# talked about emissions "per inhabitant of the territory" which was wrong;
# replace it with the corrected text about emissions per 1 sq.km of territory.
$ws.Range("B12").Value = "Выбросы загрязняющих веществ в атмосферный воздух от стационарных источников, в расчете на 1 кв.км территории – это общий объем загрязняющих выбросов в атмосферу стационарными источниками в региональном разрезе в расчете на 1 кв.км."

# Move the active selection from B2 to A2.
$ws.Range("A2").Select()
